$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fuzzy-match results in row 2 with new matched values
$ws.Range("A2").Value = "MAYFAIR ELEMENTARY SCHOOL"
$ws.Range("C2").Value = "2901 PRINCETON AVE"
$ws.Range("E2").Value = 0.73

# Re-apply the header font so the font-family id gets recorded (family=2 / Swiss for Calibri)
$ws.Range("A1:F1").Font.Name = "Calibri"

# Columns now need to accommodate the new (longer) values -> autofit like Excel did
$ws.Columns("A:F").AutoFit()
